$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the "we [can not] specify the" runs (with proofErr spell-
# check wrapping around "can not") into a single run of text, dropping the
# now-unneeded proofErr tags. The whole containing list-paragraph is located
# via Find and then replaced wholesale (via InsertXML) with an identical copy
# except for that merge, so every other run in the paragraph is left exactly
# as it was.
# ---------------------------------------------------------------------------
$search1 = $d.Content.Duplicate
$null = $search1.Find.Execute("When we use the directional light, spot light and other light sources, we")
$para1Start = $search1.Paragraphs(1).Range.Duplicate
$para1Start.Collapse(1)

$para1Xml = @'
<w:p w:rsidR="00027EBA" w:rsidRDefault="0045044C" w:rsidP="0045044C"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="42"/></w:numPr><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">When we use the directional light, spot light and other light sources, we can not specify the </w:t></w:r><w:r w:rsidRPr="00277694"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="FF0000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Spectral </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="FF0000"/><w:sz w:val="24"/></w:rPr><w:t>curve</w:t></w:r><w:r w:rsidRPr="0045044C"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">but we can measure the absorption of different wavelength as shown in the daylight example. However, I was wondering what </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00027EBA"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>are the default spectral curves</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> for those light sources</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">.  </w:t></w:r><w:r w:rsidR="00027EBA"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>U</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>niform distribution for all wavelengths or n</w:t></w:r><w:r w:rsidR="00027EBA"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">ot? </w:t></w:r></w:p>
'@

$para1Start.InsertXML($para1Xml)

# ---------------------------------------------------------------------------
# Change 2: insert four new log paragraphs right before the existing
# paragraph that carries the "_GoBack" bookmark, relocate the bookmark into a
# brand-new paragraph of its own, and leave one new empty paragraph behind.
# Because InsertXML replaces the whole paragraph touched by a collapsed
# range, the fragment below starts by re-creating the original bookmark
# paragraph (same rsid attributes) stripped down to just its <w:pPr>, exactly
# as the target diff shows.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bmParaStart = $bm.Range.Paragraphs(1).Range.Duplicate
$bmParaStart.Collapse(1)

$fragment2Xml = @'
<w:p w:rsidR="00F511C1" w:rsidRDefault="00F511C1" w:rsidP="008A690A"><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">2016-4-19 due to the differences in versions, e.g. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>waterflux</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> optimization, carbon optimization, visualization, I decide to create a latest one has all the updated functions. Since the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>waterflux</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> optimization works well and close to end, so I do not update this.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:t>Furthermore, I created a Boolean variable to control the running differences when read external data or calculate environment variables the model itself.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:t>2016-4-19</w:t></w:r><w:r><w:t xml:space="preserve"> I updated the input file to make everything combines into one. </w:t></w:r><w:r><w:t>I also add lots of running controls in the model input file. Note calculation equations are not readable</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr></w:p>
'@

$bmParaStart.InsertXML($fragment2Xml)

